$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user_data")

# Fix the admin username value (was "Adminsss", should be "Admin")
$ws.Range("A2").Value = "Admin"

# Move the active selection to F12 (was C9)
$ws.Range("F12").Select()
